# crawford_brandon.xlsx: add a "Save" stat column (H) next to "sum" (G)
# - H1 header "Save", styled like the other header cells (copy G1's format)
# - H2 data value 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing header's formatting (bold, centered, bordered) onto H1
# by copy/paste-special of formats, then set the values.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
